$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit rotates product listing rows 5-8: row 8's data moves up to row 5,
# and the old rows 5-7 shift down to rows 6-8 (rows 9+ stay untouched).

# Row 5 <- old row 8's data (rating/E column stays blank, unchanged)
$ws.Range("A5").Value = "100% all'ingrosso dell'oem PP del panno premio di grandi dimensioni del rotolo perforato per fusione del tergicristallo industriale resistente che pulisce l'olio con il MOQ basso"
$ws.Range("B5").Value = "9,53-10,39 €"
$ws.Range("C5").Value = "Ordine minimo: 50 rulli"
$ws.Range("D5").Value = "Suzhou Minghuiyuan Technology Co., Ltd."

# Row 6 <- old row 5's data; rating (E) becomes blank (old row 5 had none)
$ws.Range("A6").Value = "Wholesale Efficient Cleaning Wiping Workshop Shop Roll Jumbo Towel Industrial Wipes"
$ws.Range("B6").Value = "12,80-13,47 €"
$ws.Range("C6").Value = "Ordine minimo: 500 parti"
$ws.Range("D6").Value = "Huzhou Auline Sanitary Material Co., Ltd."
$ws.Range("E6").ClearContents()

# Row 7 <- old row 6's data; company/rating (D/E) unchanged (same values already)
$ws.Range("A7").Value = "KILINE 30% sconto fabbrica all'ingrosso usa e getta Micro fibra grande panno per pulire olio industriale cucina senza polvere"
$ws.Range("B7").Value = "1,04 €"
$ws.Range("C7").Value = "Ordine minimo: 10 parti"

# Row 8 <- old row 7's data; company (D) unchanged, rating (E) now "5.0"
$ws.Range("A8").Value = "KILINE 30% Discount Wholesale Large Roll Dust-Free Paper Industrial Multi-Function Wiping Machine Workshop Cleaning Dust Oil"
$ws.Range("B8").Value = "0,0867-0,4851 €"
$ws.Range("C8").Value = "Ordine minimo: 10 rulli"
$ws.Range("D8").Value = "Shanghai Kiline Paper Co., Ltd."
$ws.Range("E8").Value = "'5.0"
